# Apply the May 9th changes: remove the first 3 data rows (old rows 2-4)
# and append 13 new data rows at the end of the data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 3 data rows (rows 2, 3, 4), shifting remaining data up.
$ws.Range("A2:C4").Delete()

# New rows to append after the shifted data (which now ends at row 18).
$newData = @(
    @(-5.452502425124007, 1.563533739345819, 3.775565188105511),
    @(-4.596750730421486, 7.476189779072271, 4.383058036245957),
    @(-3.8169967372243, 9.492706188341467, -5.739285922631989),
    @(5.77226390198965, -3.99830269232031, -2.008213531680217),
    @(2.478664084178635, -3.857492981887443, -1.977369552705346),
    @(0.5908174134972555, -0.1589468512229657, 2.452448280846219),
    @(-13.21499349467632, -7.594239313064515, 6.02689108034458),
    @(4.174701481330613, -10.92443460371436, 1.552455660046618),
    @(1.570928898958005, -8.095582741062824, 0.9565766177526321),
    @(8.965150342112763, 11.8617115107978, 1.040716253402757),
    @(0.2588674527842816, 4.578985643822968, 1.977885912104351),
    @(-2.832312445815016, 0.201128745497, 1.380974077596895),
    @(-4.569073584021599, -4.38156421046433, -0.2334796684544305)
)

$startRow = 19
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
